$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.961.04'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '2.682.71'
$ws.Range('E3').Value = '  -2.86%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'549.04"
$ws.Range('E5').Value = '  -4.89%  '
$ws.Range('D6').Value = "'157.48"
$ws.Range('E6').Value = '  -1.92%  '
$ws.Range('D7').Value = "'1.00"
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('E9').Value = '  -4.58%  '
$ws.Range('E10').Value = '  -2.56%  '
$ws.Range('E11').Value = '  -4.86%  '
$ws.Range('D12').Value = "'5.10"
$ws.Range('E12').Value = '  -12.90%  '
$ws.Range('D13').Value = '3.158.02'
$ws.Range('E13').Value = '  -2.86%  '
$ws.Range('D14').Value = "'25.99"
$ws.Range('E14').Value = '  -5.18%  '
$ws.Range('D15').Value = '62.790.15'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('E16').Value = '  -4.10%  '
$ws.Range('D17').Value = '2.684.38'
$ws.Range('E17').Value = '  -3.08%  '
$ws.Range('D18').Value = "'11.92"
$ws.Range('E18').Value = '  -2.17%  '
$ws.Range('D19').Value = "'4.57"
$ws.Range('E19').Value = '  -5.87%  '
$ws.Range('D20').Value = "'342.69"
$ws.Range('E20').Value = '  -4.41%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').Value = "'0.503"
$ws.Range('E23').Value = '  -5.00%  '
$ws.Range('D24').Value = "'63.28"
$ws.Range('E24').Value = '  -3.01%  '
$ws.Range('E25').Value = '  -2.24%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = "'8.15"
$ws.Range('E27').Value = '  -5.51%  '
$ws.Range('E28').Value = '  -8.06%  '
$ws.Range('E29').Value = '  -2.52%  '
$ws.Range('E30').Value = '  -3.78%  '
$ws.Range('E31').Value = '  -5.19%  '
$ws.Range('D32').Value = "'165.23"
$ws.Range('E32').Value = '  -1.97%  '
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('D34').Value = "'4.79"
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').Value = "'19.51"
$ws.Range('E35').Value = '  -3.47%  '
$ws.Range('D36').Value = "'1.42"
$ws.Range('E36').Value = '  -6.32%  '
$ws.Range('E37').Value = '  -3.75%  '
$ws.Range('D38').Value = "'338.77"
$ws.Range('E38').Value = '  -3.42%  '
$ws.Range('E39').Value = '  -3.55%  '
$ws.Range('E40').Value = '  -7.60%  '
$ws.Range('D41').Value = "'38.13"
$ws.Range('E41').Value = '  -2.16%  '
$ws.Range('E42').Value = '  -6.32%  '
$ws.Range('D43').Value = "'20.33"
$ws.Range('E43').Value = '  -5.77%  '
$ws.Range('D44').Value = "'20.72"
$ws.Range('E44').Value = '  -8.11%  '
$ws.Range('D45').Value = "'0.617"
$ws.Range('E45').Value = '  -2.22%  '
$ws.Range('D46').Value = "'0.0559"
$ws.Range('E46').Value = '  -6.12%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = "'11.04"
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').Value = "'0.0973"
$ws.Range('E49').Value = '  -3.90%  '
$ws.Range('D50').Value = "'129.35"
$ws.Range('E50').Value = '  -5.51%  '
$ws.Range('D51').Value = '2.086.80'
$ws.Range('E51').Value = '  -2.86%  '
